# Update the "Macro_taxonomy" sheet:
#  - insert a new row before row 5 (shifts the old rows 5-11 down to 6-12)
#  - populate the new row 5 with the "Block/Stone/Brick" / "Urban" / "CR/LFINF" / 0.125 record
#  - update the macro_proportion values of the existing Urban rows (2-4) to match
#    the new proportions implied by the extra row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Macro_taxonomy")

# Insert a new row at position 5 - this shifts rows 5..11 down to 6..12
$ws.Rows.Item(5).Insert()

# Recalculated macro_proportion values for the existing "Block/Stone/Brick" / "Urban" rows
$ws.Range("D2").Value = 0.3125
$ws.Range("D3").Value = 0.2500000000000001
$ws.Range("D4").Value = 0.3125

# New row 5: Block/Stone/Brick, Urban, CR/LFINF, 0.125
$ws.Range("A5").Value = "Block/Stone/Brick"
$ws.Range("B5").Value = "Urban"
$ws.Range("C5").Value = "CR/LFINF"
$ws.Range("D5").Value = 0.125
